# Insert two new weekly price rows right before the existing row 530,
# pushing all subsequent rows down by two (old 530-627 -> new 532-629).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 530 (shifts row 530..627 down to 532..629)
$ws.Range("A530:R531").EntireRow.Insert()

# New row 530 data
$ws.Range("A530").Value2 = 9
$ws.Range("B530").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C530").Value2 = "Metropolitana"
$ws.Range("D530").Value2 = 45218
$ws.Range("E530").Value2 = 13
$ws.Range("F530").Value2 = 100112032
$ws.Range("G530").Value2 = "Zapallo italiano"
$ws.Range("H530").Value2 = "Sin especificar"
$ws.Range("I530").Value2 = "Primera"
$ws.Range("J530").Value2 = 97
$ws.Range("K530").Value2 = 15000
$ws.Range("L530").Value2 = 16000
$ws.Range("M530").Value2 = 15485
$ws.Range("N530").Value2 = '$/caja 50 unidades'
$ws.Range("O530").Value2 = "Región de Arica y Parinacota"
$ws.Range("P530").Value2 = 310
$ws.Range("Q530").Value2 = 50
$ws.Range("R530").Value2 = "Hortaliza"

# New row 531 data
$ws.Range("A531").Value2 = 9
$ws.Range("B531").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C531").Value2 = "Metropolitana"
$ws.Range("D531").Value2 = 45218
$ws.Range("E531").Value2 = 13
$ws.Range("F531").Value2 = 100112032
$ws.Range("G531").Value2 = "Zapallo italiano"
$ws.Range("H531").Value2 = "Sin especificar"
$ws.Range("I531").Value2 = "Primera"
$ws.Range("J531").Value2 = 70
$ws.Range("K531").Value2 = 15000
$ws.Range("L531").Value2 = 17000
$ws.Range("M531").Value2 = 16000
$ws.Range("N531").Value2 = '$/caja 60 unidades'
$ws.Range("O531").Value2 = "Región de O'Higgins"
$ws.Range("P531").Value2 = 267
$ws.Range("Q531").Value2 = 60
$ws.Range("R531").Value2 = "Hortaliza"
